# The commit inserts one new weekly price record (a new row 63) into the
# "Feria Lagunitas de Puerto Montt - Albahaca" price log, pushing every
# existing record from row 63 onward down by one row (63->64, ..., 155->156).
#
# We reproduce this by inserting a blank row at position 63 (which shifts
# all subsequent rows down, exactly like Excel's own Insert Row command),
# and then populating that new row with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 63, shifting rows 63:155 down to 64:156.
$ws.Rows(63).Insert()

# Fill in the new record in row 63.
$ws.Range("A63").Value = 4
$ws.Range("B63").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C63").Value = "Los Lagos"
$ws.Range("D63").Value = 44915
$ws.Range("E63").Value = 10
$ws.Range("F63").Value = 100112052
$ws.Range("G63").Value = "Albahaca"
$ws.Range("H63").Value = "Sin especificar"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 90
$ws.Range("K63").Value = 7000
$ws.Range("L63").Value = 8000
$ws.Range("M63").Value = 7500
$ws.Range("N63").Value = '$/docena de matas'
$ws.Range("O63").Value = "Región Metropolitana"
$ws.Range("P63").Value = 1250
$ws.Range("Q63").Value = 6
$ws.Range("R63").Value = "Hortaliza"
